{"js": "// Update the worksheet date and every multiplication problem in the\n// table to the newly generated values. Each old value is unique in\n// the document, so a literal (non-wildcard) search-and-replace on\n// each pair, applied in document order, reproduces the diff exactly.\nconst replacements = [\n  [\"2025-07-01 Tuesday\", \"2025-07-02 Wednesday\"],\n  [\"69\u00d740=\", \"14\u00d721=\"],\n  [\"77\u00d732=\", \"58\u00d788=\"],\n  [\"86\u00d724=\", \"21\u00d764=\"],\n  [\"87\u00d712=\", \"53\u00d784=\"],\n  [\"53\u00d781=\", \"17\u00d748=\"],\n  [\"66\u00d788=\", \"51\u00d742=\"],\n  [\"67\u00d743=\", \"16\u00d788=\"],\n  [\"58\u00d794=\", \"23\u00d738=\"],\n  [\"59\u00d771=\", \"38\u00d793=\"],\n  [\"31\u00d768=\", \"65\u00d739=\"],\n  [\"16\u00d736=\", \"88\u00d770=\"],\n  [\"18\u00d744=\", \"32\u00d718=\"],\n  [\"27\u00d725=\", \"14\u00d725=\"],\n  [\"82\u00d731=\", \"29\u00d798=\"],\n  [\"46\u00d785=\", \"14\u00d752=\"],\n  [\"75\u00d717=\", \"91\u00d731=\"],\n  [\"66\u00d771=\", \"17\u00d715=\"],\n  [\"81\u00d771=\", \"31\u00d745=\"],\n  [\"20\u00d752=\", \"67\u00d756=\"],\n  [\"27\u00d749=\", \"95\u00d747=\"],\n  [\"56\u00d766=\", \"36\u00d788=\"],\n  [\"83\u00d787=\", \"58\u00d764=\"],\n  [\"33\u00d716=\", \"33\u00d712=\"],\n  [\"50\u00d738=\", \"66\u00d772=\"],\n  [\"72\u00d736=\", \"60\u00d790=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and all 25 multiplication problems in the\n# table to the newly generated values, using Find/Replace on each\n# unique \"old text\" -> \"new text\" pair (exact literal match, no\n# wildcards) so each cell is updated independently and in place.\n$d = $word.ActiveDocument\n\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"2025-07-01 Tuesday\", $true, $false, $false, $false, $false, $true, 1, $false, \"2025-07-02 Wednesday\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"69\u00d740=\", $true, $false, $false, $false, $false, $true, 1, $false, \"14\u00d721=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"77\u00d732=\", $true, $false, $false, $false, $false, $true, 1, $false, \"58\u00d788=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"86\u00d724=\", $true, $false, $false, $false, $false, $true, 1, $false, \"21\u00d764=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"87\u00d712=\", $true, $false, $false, $false, $false, $true, 1, $false, \"53\u00d784=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"53\u00d781=\", $true, $false, $false, $false, $false, $true, 1, $false, \"17\u00d748=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"66\u00d788=\", $true, $false, $false, $false, $false, $true, 1, $false, \"51\u00d742=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"67\u00d743=\", $true, $false, $false, $false, $false, $true, 1, $false, \"16\u00d788=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"58\u00d794=\", $true, $false, $false, $false, $false, $true, 1, $false, \"23\u00d738=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"59\u00d771=\", $true, $false, $false, $false, $false, $true, 1, $false, \"38\u00d793=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"31\u00d768=\", $true, $false, $false, $false, $false, $true, 1, $false, \"65\u00d739=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"16\u00d736=\", $true, $false, $false, $false, $false, $true, 1, $false, \"88\u00d770=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"18\u00d744=\", $true, $false, $false, $false, $false, $true, 1, $false, \"32\u00d718=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"27\u00d725=\", $true, $false, $false, $false, $false, $true, 1, $false, \"14\u00d725=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"82\u00d731=\", $true, $false, $false, $false, $false, $true, 1, $false, \"29\u00d798=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"46\u00d785=\", $true, $false, $false, $false, $false, $true, 1, $false, \"14\u00d752=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"75\u00d717=\", $true, $false, $false, $false, $false, $true, 1, $false, \"91\u00d731=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"66\u00d771=\", $true, $false, $false, $false, $false, $true, 1, $false, \"17\u00d715=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"81\u00d771=\", $true, $false, $false, $false, $false, $true, 1, $false, \"31\u00d745=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"20\u00d752=\", $true, $false, $false, $false, $false, $true, 1, $false, \"67\u00d756=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"27\u00d749=\", $true, $false, $false, $false, $false, $true, 1, $false, \"95\u00d747=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"56\u00d766=\", $true, $false, $false, $false, $false, $true, 1, $false, \"36\u00d788=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"83\u00d787=\", $true, $false, $false, $false, $false, $true, 1, $false, \"58\u00d764=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"33\u00d716=\", $true, $false, $false, $false, $false, $true, 1, $false, \"33\u00d712=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"50\u00d738=\", $true, $false, $false, $false, $false, $true, 1, $false, \"66\u00d772=\", 2) | Out-Null\n$r = $d.Content\n$r.Find.ClearFormatting()\n$r.Find.Execute(\"72\u00d736=\", $true, $false, $false, $false, $false, $true, 1, $false, \"60\u00d790=\", 2) | Out-Null\n"}
